$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.496.02"
$ws.Range("E2").Value = "  -2.00%  "

# Row 3
$ws.Range("D3").Value = "1.749.00"
$ws.Range("E3").Value = "  -2.23%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.38"
$ws.Range("E5").Value = "  +0.29%  "

# Row 6
$ws.Range("E6").Value = "  -0.05%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4455"
$ws.Range("E7").Value = "  +3.71%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3601"
$ws.Range("E8").Value = "  -0.42%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07497"
$ws.Range("E9").Value = "  -0.02%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.95"
$ws.Range("E10").Value = "  -6.28%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.090"
$ws.Range("E11").Value = "  -2.10%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.04%  "

# Row 13
$ws.Range("E13").Value = "  -4.86%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.022"
$ws.Range("E14").Value = "  -2.00%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.116"
$ws.Range("E15").Value = "  -2.68%  "

# Row 16
$ws.Range("D16").Value = "1.749.96"
$ws.Range("E16").Value = "  -2.31%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.29"
$ws.Range("E17").Value = "  +1.18%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001060"
$ws.Range("E18").Value = "  -0.48%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06396"
$ws.Range("E19").Value = "  +0.61%  "

# Row 20
$ws.Range("E20").Value = "  +0.00%  "

# Row 21
$ws.Range("E21").Value = "  -2.66%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.849"
$ws.Range("E22").Value = "  -2.07%  "

# Row 23
$ws.Range("D23").Value = "27.551.53"
$ws.Range("E23").Value = "  -1.87%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.16"
$ws.Range("E24").Value = "  -1.90%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.110"
$ws.Range("E25").Value = "  -0.88%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.35"
$ws.Range("E26").Value = "  +1.49%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.47"
$ws.Range("E27").Value = "  +0.56%  "

# Row 28
$ws.Range("D28").Value = "1.950.64"
$ws.Range("E28").Value = "  -2.31%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.092"
$ws.Range("E29").Value = "  -3.66%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.14"
$ws.Range("E30").Value = "  -1.64%  "

# Row 31
$ws.Range("E31").Value = "  -6.60%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.656"
$ws.Range("E32").Value = "  +3.84%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09012"
$ws.Range("E33").Value = "  +0.32%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.535"
$ws.Range("E34").Value = "  -3.81%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.97"
$ws.Range("E35").Value = "  -5.13%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02294"
$ws.Range("E36").Value = "  -1.23%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06014"
$ws.Range("E37").Value = "  -0.63%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2083"
$ws.Range("E38").Value = "  -1.27%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6335"
$ws.Range("E39").Value = "  -1.92%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.946"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.209"
$ws.Range("E41").Value = "  +1.81%  "

# Row 42
$ws.Range("B42").Value = "Frax"
$ws.Range("C42").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  +0.02%  "

# Row 43
$ws.Range("B43").Value = "WEMIXTOKEN"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.382"
$ws.Range("E43").Value = "  -2.66%  "

# Row 44
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.762"
$ws.Range("E44").Value = "  -1.06%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.13"
$ws.Range("E45").Value = "  -3.30%  "

# Row 46
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.712"
$ws.Range("E46").Value = "  +0.32%  "

# Row 47
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5878"
$ws.Range("E47").Value = "  -1.92%  "

# Row 48
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "121.87"
$ws.Range("E48").Value = "  -2.12%  "

# Row 49
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.951"
$ws.Range("E49").Value = "  -1.46%  "

# Row 50
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.146"
$ws.Range("E50").Value = "  -0.73%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06848"
$ws.Range("E51").Value = "  -1.48%  "
